$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "sub-s1"
$ws.Range("A3").Value = "sub-s2"
$ws.Range("A4").Select()
